$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column G: "New COVID Hospitalizations" ---
$ws.Cells.Item(1, 7).Value = "New COVID Hospitalizations"

# Copy the header style (bold) from F1 to G1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set column G values for existing rows (6-28), mirroring data column F's number style
$gValues = @{
    6  = 489
    7  = 823
    8  = 776
    9  = 909
    10 = 1084
    11 = 1796
    12 = 1813
    13 = 1722
    14 = 2241
    15 = 1883
    16 = 2507
    17 = 2844
    18 = 2857
    19 = 3413
    20 = 3261
    21 = 2821
    22 = 2082
    23 = 2553
    24 = 3034
    25 = 2848
    26 = 2882
    27 = 2486
    28 = 2538
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
    $ws.Range("F$row").Copy()
    $ws.Range("G$row").PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# --- Add new row 29 with hospitalization data for 2020-04-12 ---
# Copy formats per-column from row 28 (skip column E, which stays empty on row 29)
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B28").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C28").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D28").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F28").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G28").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A29").Value = 43933
$ws.Range("B29").Value = 118
$ws.Range("C29").Value = -42
$ws.Range("D29").Value = -21
$ws.Range("F29").Value = 671
$ws.Range("G29").Value = 1958

# --- Adjust column width for new column G ---
$ws.Columns.Item(7).ColumnWidth = 24.73

# --- Update selection to reflect new active cell ---
$ws.Range("G30").Select()
